# 1) Rename the "TeachingStuff" department to "Teaching Staff" and add a new
#    "Student Affairs" department to the Departments lookup sheet.
# 2) Record a sample/test row on Sheet1 using the new "Student Affairs" value
#    (role condition / department sample data for the employee index import).

$wb  = $excel.ActiveWorkbook
$deptWs = $wb.Worksheets.Item("Departments")

# Rename existing department label.
$deptWs.Range("A1").Value = "Teaching Staff"

# Add the new department as a new row, copying formatting from the row above.
$deptWs.Range("A4").Copy($deptWs.Range("A5"))
$deptWs.Range("A5").Value = "Student Affairs"

# Add a sample row on Sheet1 referencing the new department.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Range("E4").Value = "Student Affairs"

# Make Sheet1 active with E4 selected, matching the saved selection state.
$sheet1.Activate()
[void]$sheet1.Range("E4").Select()
